# Add the "2nd" function row (row 26) to the calculator operations sheet,
# mirroring the formatting already used for the rows above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (font/alignment/border/number-format) from the last
# existing data row (25) onto the new row's two populated cells, A26 and
# D26, without disturbing B26/C26 (which stay empty, just like the diff).
$ws.Range("A25").Copy() | Out-Null
$ws.Range("A26").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("D25").Copy() | Out-Null
$ws.Range("D26").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Fill in the new row's values.
$ws.Range("A26").Value = "2nd"
$ws.Range("D26").Value = "OPENS NEW SCREEN"

# Match the row height used by the rest of the table.
$ws.Rows("26").RowHeight = 15.75

# The active selection in the saved file moved from D28 to D27.
$ws.Range("D27").Select() | Out-Null
